# Remove the redundancy of the code
$wb = $excel.ActiveWorkbook

# --- RateCard & Pricing: move the selection from D5 to E1 ---
$wsRate = $wb.Worksheets.Item("RateCard & Pricing")
$wsRate.Activate()
$wsRate.Range("E1").Select()

# --- Customer: selection stays at D2, sheet is no longer the active tab ---
$wsCustomer = $wb.Worksheets.Item("Customer")
$wsCustomer.Activate()
$wsCustomer.Range("D2").Select()

# --- CreateCategory: the "Plan Pricing" group (I1:M1) duplicated E1's
#     "Test Product" in L1. Fold the redundant single-row layout into a
#     second row and drop the duplicate Test Product cell. ---
$wsCat = $wb.Worksheets.Item("CreateCategory")

# Plan Pricing Category (was I1) -> D2
$wsCat.Range("I1").Copy()
$wsCat.Range("D2").PasteSpecial()

# Plan Pricing Product (was J1) -> E2
$wsCat.Range("J1").Copy()
$wsCat.Range("E2").PasteSpecial()

# "20" (was K1) -> G2
$wsCat.Range("K1").Copy()
$wsCat.Range("G2").PasteSpecial()

# Plan Pricing (was M1) -> I1 (replacing "Plan Pricing Category")
$wsCat.Range("M1").Copy()
$wsCat.Range("I1").PasteSpecial()

# L1 was a redundant duplicate of E1 ("Test Product") - drop J1:M1 entirely
$wsCat.Range("J1:M1").Clear()

# CreateCategory becomes the active sheet/tab, selection at I1
$wsCat.Activate()
$wsCat.Range("I1").Select()
